# Se arregló la apariencia de todos los submenús:
# apply yellow highlighting to the remaining "Requisitos funcionales del
# sistema" bullets so every sub-item in that list looks consistent with
# its already-highlighted siblings (both the run text and the paragraph
# mark pick up the highlight, matching how Word records it when the whole
# line, pilcrow included, is highlighted).

$d = $word.ActiveDocument

$targets = @(
    "Matricular estudiantes en cursos específicos.",
    "Consultar los estudiantes inscritos en un curso.",
    "Consultar los horarios asignados a un curso.",
    "Eliminar matrículas sin afectar la información de los estudiantes ni de los cursos."
)

$runPattern = [regex]'(?s)<w:r(\s[^>]*)?>(.*?)</w:r>'
$rprPattern = [regex]'(?s)^<w:rPr>(.*?)</w:rPr>'

foreach ($text in $targets) {
    $searchRng = $d.Content
    $found = $searchRng.Find.Execute($text, $true, $false, $false, $false,
                                      $false, $true, 1, $false, "", 0)
    if (-not $found) {
        continue
    }
    $hitStart = $searchRng.Start

    # Resolve the (unclipped) paragraph that contains the hit so its
    # Range spans the whole paragraph, including the trailing paragraph
    # mark -- a range narrowed by Find.Execute does not.
    $allParas = $d.Paragraphs
    $rng = $null
    for ($i = 1; $i -le $allParas.Count; $i++) {
        $p = $allParas.Item($i)
        if ($p.Range.Start -le $hitStart -and $p.Range.End -gt $hitStart) {
            $rng = $p.Range
            break
        }
    }
    if ($rng -eq $null) {
        continue
    }

    $full = $rng.WordOpenXML
    if (-not ($full -match '(?s)(<w:p [^>]*>.*?</w:p>)')) {
        continue
    }
    $pxml = $matches[1]

    # Make sure the paragraph mark itself (w:pPr/w:rPr) carries the
    # highlight, creating the w:rPr element if the paragraph didn't have
    # one yet.
    if ($pxml -match '(?s)<w:pPr>(.*?)</w:pPr>') {
        $pprInner = $matches[1]
        if ($pprInner -match '(?s)<w:rPr>(.*?)</w:rPr>') {
            $newPprInner = $pprInner -replace '(?s)<w:rPr>(.*?)</w:rPr>', '<w:rPr>$1<w:highlight w:val="yellow"/></w:rPr>'
        } else {
            $newPprInner = $pprInner + '<w:rPr><w:highlight w:val="yellow"/></w:rPr>'
        }
        $pxml = $pxml -replace '(?s)<w:pPr>.*?</w:pPr>', ('<w:pPr>' + $newPprInner + '</w:pPr>')
    } else {
        $pxml = $pxml -replace '(?s)(<w:p [^>]*>)', ('$1<w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>')
    }

    # Highlight every run in the paragraph, creating w:rPr when missing.
    $ms = $runPattern.Matches($pxml)
    $rebuilt = ""
    $lastEnd = 0
    foreach ($m in $ms) {
        $rebuilt += $pxml.Substring($lastEnd, $m.Index - $lastEnd)
        $attrs = $m.Groups[1].Value
        $inner = $m.Groups[2].Value
        $rm = $rprPattern.Match($inner)
        if ($rm.Success) {
            $rprInner = $rm.Groups[1].Value
            $newInner = $rprPattern.Replace($inner, ('<w:rPr>' + $rprInner + '<w:highlight w:val="yellow"/></w:rPr>'), 1)
        } else {
            $newInner = '<w:rPr><w:highlight w:val="yellow"/></w:rPr>' + $inner
        }
        $rebuilt += "<w:r$attrs>$newInner</w:r>"
        $lastEnd = $m.Index + $m.Length
    }
    $rebuilt += $pxml.Substring($lastEnd)

    $rng.InsertXML($rebuilt)
}
